# Auto-generated edit script: update crypto price/volume snapshot values.
# Matches the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    # Force text interpretation so numeric-looking strings (e.g. "1.011")
    # are not silently coerced into numbers by Excel's input parser.
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = $origStyle
}

$ws.Range("D2").Value = "30.632.57"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "2.117.48"
$ws.Range("E3").Value = "  +1.18%  "
Set-TextValue "D4" "1.011"
$ws.Range("E4").Value = "  +0.78%  "
Set-TextValue "D5" "336.85"
$ws.Range("E5").Value = "  +2.09%  "
Set-TextValue "D6" "1.010"
$ws.Range("E6").Value = "  +0.75%  "
$ws.Range("E7").Value = "  +0.64%  "
Set-TextValue "D8" "0.4548"
$ws.Range("E8").Value = "  +3.06%  "
Set-TextValue "D9" "54.48"
$ws.Range("E9").Value = "  +0.54%  "
Set-TextValue "D10" "0.09123"
$ws.Range("E10").Value = "  +2.11%  "
$ws.Range("E11").Value = "  +2.06%  "
Set-TextValue "D12" "24.47"
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("D13").Value = "2.116.53"
$ws.Range("E13").Value = "  +0.92%  "
Set-TextValue "D14" "6.858"
$ws.Range("E14").Value = "  +2.67%  "
Set-TextValue "D15" "8.146"
$ws.Range("E15").Value = "  +5.81%  "
$ws.Range("E16").Value = "  +4.97%  "
Set-TextValue "D17" "97.14"
$ws.Range("E17").Value = "  +1.29%  "
Set-TextValue "D18" "1.011"
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("E19").Value = "  +0.98%  "
$ws.Range("E20").Value = "  +1.57%  "
$ws.Range("E21").Value = "  +0.68%  "
Set-TextValue "D22" "6.306"
$ws.Range("E22").Value = "  +0.77%  "
$ws.Range("D23").Value = "30.690.64"
$ws.Range("E23").Value = "  +0.59%  "
Set-TextValue "D24" "12.94"
$ws.Range("E24").Value = "  +5.04%  "
Set-TextValue "D25" "2.357"
$ws.Range("E25").Value = "  +2.12%  "
$ws.Range("D26").Value = "2.365.75"
$ws.Range("E26").Value = "  +1.05%  "
$ws.Range("E27").Value = "  +0.98%  "
Set-TextValue "D28" "164.41"
$ws.Range("E28").Value = "  +0.38%  "
Set-TextValue "D29" "2.553"
$ws.Range("E29").Value = "  -0.55%  "
Set-TextValue "D30" "134.67"
Set-TextValue "D31" "1.212"
$ws.Range("E31").Value = "  +1.83%  "
$ws.Range("E32").Value = "  +0.36%  "
Set-TextValue "D33" "1.647"
$ws.Range("E33").Value = "  +0.01%  "
Set-TextValue "D34" "6.361"
$ws.Range("E34").Value = "  +3.41%  "
Set-TextValue "D35" "3.943"
$ws.Range("E35").Value = "  +1.14%  "
$ws.Range("E36").Value = "  +5.28%  "
Set-TextValue "D37" "5.889"
$ws.Range("E37").Value = "  +7.75%  "
Set-TextValue "D38" "0.02629"
$ws.Range("E38").Value = "  +2.84%  "
Set-TextValue "D39" "0.06848"
$ws.Range("E39").Value = "  +0.66%  "
$ws.Range("E40").Value = "  +3.27%  "
Set-TextValue "D41" "12.61"
$ws.Range("E41").Value = "  -0.01%  "
Set-TextValue "D42" "0.6895"
$ws.Range("E42").Value = "  +0.42%  "
$ws.Range("E43").Value = "  +0.65%  "
Set-TextValue "D44" "14.84"
$ws.Range("E44").Value = "  +6.25%  "
Set-TextValue "D45" "0.6491"
$ws.Range("E45").Value = "  +2.68%  "
Set-TextValue "D46" "2.315"
$ws.Range("E46").Value = "  +5.53%  "
$ws.Range("E47").Value = "  +20.80%  "
Set-TextValue "D48" "3.693"
$ws.Range("E48").Value = "  +1.76%  "
Set-TextValue "D49" "1.256"
$ws.Range("E49").Value = "  +0.88%  "
Set-TextValue "D50" "83.39"
$ws.Range("E50").Value = "  +2.24%  "
$ws.Range("E51").Value = "  -3.80%  "
